# "Module / memory-profiler" slide (internal sldId 799, creationId
# 2357954353) — nudge its picture ("Picture 6", shape id 7) up slightly:
# off/y moves from 367083 EMU to 353333 EMU (off/x, ext stay untouched).
#
# PowerPoint's object model reports shape position/size in points, while
# the OOXML stores EMU (1 pt = 12700 EMU), so convert before assigning.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(31)

$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Picture 6") {
        $pic = $shp
        break
    }
}
if ($pic -eq $null) {
    $pic = $s.Shapes.Item("Picture 6")
}

$targetTopEmu = 353333
$pic.Top = $targetTopEmu / 12700
